$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "568.81" would otherwise be auto-coerced to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '63.642.71'
$ws.Range("E2").Value = '  +1.22%  '

# Row 3
$ws.Range("D3").Value = '3.424.56'
$ws.Range("E3").Value = '  +2.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '570.96'
$ws.Range("E5").Value = '  +2.67%  '

# Row 6
$ws.Range("D6").Value = '155.97'
$ws.Range("E6").Value = '  +2.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = '3.427.07'
$ws.Range("E8").Value = '  +2.52%  '

# Row 9
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  +2.66%  '

# Row 10
$ws.Range("D10").Value = '7.45'
$ws.Range("E10").Value = '  -0.55%  '

# Row 11
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").Value = '  +3.67%  '

# Row 12
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +0.10%  '

# Row 13
$ws.Range("D13").Value = '4.017.95'
$ws.Range("E13").Value = '  +2.75%  '

# Row 14
$ws.Range("D14").Value = '0.134'
$ws.Range("E14").Value = '  -2.94%  '

# Row 15
$ws.Range("D15").Value = '0.0000191'
$ws.Range("E15").Value = '  +4.95%  '

# Row 16
$ws.Range("D16").Value = '27.16'
$ws.Range("E16").Value = '  +0.98%  '

# Row 17
$ws.Range("D17").Value = '63.796.61'
$ws.Range("E17").Value = '  +1.50%  '

# Row 18
$ws.Range("D18").Value = '3.395.96'
$ws.Range("E18").Value = '  +2.31%  '

# Row 19
$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  -2.18%  '

# Row 20
$ws.Range("D20").Value = '14.21'
$ws.Range("E20").Value = '  +3.27%  '

# Row 21
$ws.Range("D21").Value = '386.23'
$ws.Range("E21").Value = '  -0.55%  '

# Row 22
$ws.Range("D22").Value = '8.25'
$ws.Range("E22").Value = '  -2.31%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.06%  '

# Row 24
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '0.538'
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '72.13'
$ws.Range("E25").Value = '  +2.12%  '

# Row 26
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  +23.06%  '

# Row 27
$ws.Range("D27").Value = '9.49'
$ws.Range("E27").Value = '  +7.90%  '

# Row 28
$ws.Range("E28").Value = '  -2.02%  '

# Row 29
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.20%  '

# Row 30
$ws.Range("D30").Value = '6.04'
$ws.Range("E30").Value = '  +8.51%  '

# Row 31
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '1.37'
$ws.Range("E31").Value = '  +5.15%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '2.01'
$ws.Range("E32").Value = '  +1.16%  '

# Row 33
$ws.Range("D33").Value = '23.41'
$ws.Range("E33").Value = '  +1.85%  '

# Row 34
$ws.Range("D34").Value = '6.45'
$ws.Range("E34").Value = '  +0.95%  '

# Row 35
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.08%  '

# Row 36
$ws.Range("D36").Value = '6.92'
$ws.Range("E36").Value = '  +3.43%  '

# Row 37
$ws.Range("D37").Value = '159.63'
$ws.Range("E37").Value = '  -0.57%  '

# Row 38
$ws.Range("D38").Value = '1.46'
$ws.Range("E38").Value = '  -1.82%  '

# Row 39
$ws.Range("D39").Value = '0.0771'
$ws.Range("E39").Value = '  +4.55%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '2.913.11'
$ws.Range("E40").Value = '  +3.31%  '

# Row 41
$ws.Range("D41").Value = '1.84'
$ws.Range("E41").Value = '  -2.68%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '27.06'
$ws.Range("E42").Value = '  -0.15%  '

# Row 43
$ws.Range("E43").Value = '  +1.42%  '

# Row 44
$ws.Range("D44").Value = '4.39'
$ws.Range("E44").Value = '  +1.71%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '41.56'
$ws.Range("E45").Value = '  +2.24%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.763'
$ws.Range("E46").Value = '  +2.25%  '

# Row 47
$ws.Range("D47").Value = '23.51'
$ws.Range("E47").Value = '  +7.43%  '

# Row 48
$ws.Range("D48").Value = '1.08'
$ws.Range("E48").Value = '  +3.91%  '

# Row 49
$ws.Range("D49").Value = '2.17'
$ws.Range("E49").Value = '  +20.69%  '

# Row 50
$ws.Range("E50").Value = '  +3.20%  '

# Row 51
$ws.Range("D51").Value = '0.844'
$ws.Range("E51").Value = '  +5.46%  '
